$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K") for rows 2-22
$values = @{
    2  = 4
    3  = 3
    4  = 4
    5  = 9
    6  = 3
    7  = 5
    8  = 4
    9  = 2
    10 = 5
    11 = 9
    12 = 7
    13 = 5
    14 = 3
    15 = 2
    16 = 4
    17 = 6
    18 = 7
    19 = 3
    20 = 1
    21 = 6
    22 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
